# ============================================================
# Otomatik guncelleme: 2025-06-19 02:42:11
# Refresh "durum" start dates, append a new 2025-06-19 price
# observation to each district sheet, and refresh the
# "eskalasyon" (escalation) summary sheet to include it.
# ============================================================

$wb = $excel.ActiveWorkbook

# --- durum: bump StartDate (col B) for every district row ---
$wsDurum = $wb.Worksheets.Item("durum")
$durumRows = $wsDurum.Cells.Item(1,1).CurrentRegion.Rows.Count
for ($r = 2; $r -le $durumRows; $r++) {
    $cell = $wsDurum.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = "2025-06-19"
}

# --- district price-history sheets: append the new 2025-06-19 observation ---
$ws = $wb.Worksheets.Item("934015")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 51.37

$ws = $wb.Worksheets.Item("065001")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 53.62

$ws = $wb.Worksheets.Item("035001")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 52.57

$ws = $wb.Worksheets.Item("055001")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 52.85

$ws = $wb.Worksheets.Item("021001")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 53.57

$ws = $wb.Worksheets.Item("038001")
$newRow = 195
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45827
$dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
$ws.Cells.Item($newRow, 3).Value = 53.25

# --- eskalasyon: refresh the escalation summary rows (2-36) ---
$wsE = $wb.Worksheets.Item("eskalasyon")

$r = 2
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 49.23
$wsE.Cells.Item($r, 4).Value = 0.06374243733794294
$wsE.Cells.Item($r, 5).Value = 0.06374243733794294
$wsE.Cells.Item($r, 6).Value = "Servis Diyarbakır"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 3
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 46.38
$wsE.Cells.Item($r, 4).Value = -0.0578915295551492
$wsE.Cells.Item($r, 5).Value = -0.0578915295551492
$wsE.Cells.Item($r, 6).Value = "Servis Diyarbakır"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 4
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 49.96
$wsE.Cells.Item($r, 4).Value = 0.07718844329452357
$wsE.Cells.Item($r, 5).Value = 0.07718844329452357
$wsE.Cells.Item($r, 6).Value = "Servis Diyarbakır"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 5
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 53.57
$wsE.Cells.Item($r, 4).Value = 0.07225780624499589
$wsE.Cells.Item($r, 5).Value = 0.07225780624499589
$wsE.Cells.Item($r, 6).Value = "Servis Diyarbakır"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 6
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 48.88
$wsE.Cells.Item($r, 4).Value = 0.05663640293990491
$wsE.Cells.Item($r, 5).Value = 0.05663640293990491
$wsE.Cells.Item($r, 6).Value = "Servis Kayseri"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 7
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 46.03
$wsE.Cells.Item($r, 4).Value = -0.05830605564648117
$wsE.Cells.Item($r, 5).Value = -0.05830605564648117
$wsE.Cells.Item($r, 6).Value = "Servis Kayseri"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 8
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 49.64
$wsE.Cells.Item($r, 4).Value = 0.07842711275255265
$wsE.Cells.Item($r, 5).Value = 0.07842711275255265
$wsE.Cells.Item($r, 6).Value = "Servis Kayseri"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 9
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 53.25
$wsE.Cells.Item($r, 4).Value = 0.07272360999194194
$wsE.Cells.Item($r, 5).Value = 0.07272360999194194
$wsE.Cells.Item($r, 6).Value = "Servis Kayseri"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 10
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 48.57
$wsE.Cells.Item($r, 4).Value = 0.05678851174934718
$wsE.Cells.Item($r, 5).Value = 0.05678851174934718
$wsE.Cells.Item($r, 6).Value = "Servis Samsun"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 11
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 45.66
$wsE.Cells.Item($r, 4).Value = -0.05991352686843743
$wsE.Cells.Item($r, 5).Value = -0.05991352686843743
$wsE.Cells.Item($r, 6).Value = "Servis Samsun"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 12
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 49.24
$wsE.Cells.Item($r, 4).Value = 0.07840560665790641
$wsE.Cells.Item($r, 5).Value = 0.07840560665790641
$wsE.Cells.Item($r, 6).Value = "Servis Samsun"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 13
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 52.85
$wsE.Cells.Item($r, 4).Value = 0.07331437855402112
$wsE.Cells.Item($r, 5).Value = 0.07331437855402112
$wsE.Cells.Item($r, 6).Value = "Servis Samsun"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 14
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.24
$wsE.Cells.Item($r, 4).Value = 0.07023108291798819
$wsE.Cells.Item($r, 5).Value = 0.07023108291798819
$wsE.Cells.Item($r, 6).Value = "Spot Araç Anadolu Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 15
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.26
$wsE.Cells.Item($r, 4).Value = -0.06308213378492811
$wsE.Cells.Item($r, 5).Value = -0.06308213378492811
$wsE.Cells.Item($r, 6).Value = "Spot Araç Anadolu Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 16
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.07907817442385912
$wsE.Cells.Item($r, 5).Value = 0.07907817442385912
$wsE.Cells.Item($r, 6).Value = "Spot Araç Anadolu Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 17
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "Spot Araç Anadolu Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 18
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.24
$wsE.Cells.Item($r, 4).Value = 0.07023108291798819
$wsE.Cells.Item($r, 5).Value = 0.07023108291798819
$wsE.Cells.Item($r, 6).Value = "Spot Araç Avrupa&Anadolu"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 19
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.26
$wsE.Cells.Item($r, 4).Value = -0.06308213378492811
$wsE.Cells.Item($r, 5).Value = -0.06308213378492811
$wsE.Cells.Item($r, 6).Value = "Spot Araç Avrupa&Anadolu"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 20
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.07907817442385912
$wsE.Cells.Item($r, 5).Value = 0.07907817442385912
$wsE.Cells.Item($r, 6).Value = "Spot Araç Avrupa&Anadolu"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 21
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "Spot Araç Avrupa&Anadolu"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 22
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.24
$wsE.Cells.Item($r, 4).Value = 0.07023108291798819
$wsE.Cells.Item($r, 5).Value = 0.07023108291798819
$wsE.Cells.Item($r, 6).Value = "Spot Araç Teknosa"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 23
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.26
$wsE.Cells.Item($r, 4).Value = -0.06308213378492811
$wsE.Cells.Item($r, 5).Value = -0.06308213378492811
$wsE.Cells.Item($r, 6).Value = "Spot Araç Teknosa"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 24
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.07907817442385912
$wsE.Cells.Item($r, 5).Value = 0.07907817442385912
$wsE.Cells.Item($r, 6).Value = "Spot Araç Teknosa"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 25
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "Spot Araç Teknosa"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 26
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.24
$wsE.Cells.Item($r, 4).Value = 0.07023108291798819
$wsE.Cells.Item($r, 5).Value = 0.07023108291798819
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 27
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.26
$wsE.Cells.Item($r, 4).Value = -0.06308213378492811
$wsE.Cells.Item($r, 5).Value = -0.06308213378492811
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 28
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.07907817442385912
$wsE.Cells.Item($r, 5).Value = 0.07907817442385912
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 29
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 30
$wsE.Cells.Item($r, 1).Value = 45784
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.2
$wsE.Cells.Item($r, 4).Value = -0.05089113162980452
$wsE.Cells.Item($r, 5).Value = -0.05089113162980452
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa İade Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 31
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.0805429864253393
$wsE.Cells.Item($r, 5).Value = 0.0805429864253393
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa İade Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 32
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa İade Toplama"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 33
$wsE.Cells.Item($r, 1).Value = 45673
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.24
$wsE.Cells.Item($r, 4).Value = 0.05587840858292359
$wsE.Cells.Item($r, 5).Value = 0.05587840858292359
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 34
$wsE.Cells.Item($r, 1).Value = 45756
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 44.26
$wsE.Cells.Item($r, 4).Value = -0.06308213378492811
$wsE.Cells.Item($r, 5).Value = -0.06308213378492811
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 35
$wsE.Cells.Item($r, 1).Value = 45822
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 47.76
$wsE.Cells.Item($r, 4).Value = 0.07907817442385912
$wsE.Cells.Item($r, 5).Value = 0.07907817442385912
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$wsE.Cells.Item($r, 7).Value = 0.05

$r = 36
$wsE.Cells.Item($r, 1).Value = 45827
$wsE.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsE.Cells.Item($r, 2).Value = "Motorin UltraForce"
$wsE.Cells.Item($r, 3).Value = 51.37
$wsE.Cells.Item($r, 4).Value = 0.07558626465661633
$wsE.Cells.Item($r, 5).Value = 0.07558626465661633
$wsE.Cells.Item($r, 6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$wsE.Cells.Item($r, 7).Value = 0.05

